# Update "想去人数" (want-to-go count) figures on the 展览 and 全部类型 sheets
# to match the latest scrape (commit 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibition) ---
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 5569
$wsExhibit.Range("F3").Value = 12805
$wsExhibit.Range("F4").Value = 310
$wsExhibit.Range("F5").Value = 626
$wsExhibit.Range("F6").Value = 198
$wsExhibit.Range("F7").Value = 386
$wsExhibit.Range("F8").Value = 1201

# --- Sheet "全部类型" (All Types) ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 5569
$wsAll.Range("F4").Value = 12805
$wsAll.Range("F5").Value = 310
$wsAll.Range("F6").Value = 626
$wsAll.Range("F7").Value = 198
$wsAll.Range("F10").Value = 386
$wsAll.Range("F11").Value = 1201
